# Actualización desde MV -datos-
# Appends new daily UF bond rate rows (04-08-2021 .. 02-09-2021) to Sheet1,
# continuing directly after the existing last row (147).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data extracted from the source update: each entry is the row's date (col A, stored
# as text) plus whichever of columns B-G have a value for that row.
$newData = @(
    @{ "A"="04-08-2021"; "C"="-0.89"; "D"="0.7";  "E"="1.6";  "F"="2.24" },
    @{ "A"="05-08-2021"; "C"="-0.92"; "D"="0.75"; "E"="1.62"; "F"="2.25" },
    @{ "A"="06-08-2021"; "C"="-1.19"; "D"="0.63"; "E"="1.56" },
    @{ "A"="09-08-2021"; "C"="-1.19"; "D"="0.71"; "E"="1.61"; "F"="2.27"; "G"="2.4" },
    @{ "A"="10-08-2021"; "C"="-1.19"; "D"="0.78"; "E"="1.64"; "F"="2.29" },
    @{ "A"="11-08-2021"; "C"="-1.22"; "D"="0.75"; "E"="1.61"; "F"="2.32"; "G"="2.42" },
    @{ "A"="12-08-2021"; "B"="-1.7"; "C"="-1.23"; "D"="0.68"; "E"="1.58"; "F"="2.3" },
    @{ "A"="13-08-2021"; "C"="-1.24"; "D"="0.61"; "E"="1.55"; "F"="2.2"; "G"="2.34" },
    @{ "A"="16-08-2021"; "B"="-1.8"; "C"="-1.34"; "D"="0.53"; "E"="1.49"; "F"="2.14"; "G"="2.29" },
    @{ "A"="17-08-2021"; "C"="-1.4"; "D"="0.51"; "E"="1.48" },
    @{ "A"="18-08-2021"; "C"="-1.44"; "D"="0.57"; "E"="1.5"; "F"="2.17"; "G"="2.33" },
    @{ "A"="19-08-2021"; "C"="-1.35"; "D"="0.63" },
    @{ "A"="20-08-2021"; "C"="-1.33"; "D"="0.62"; "E"="1.62"; "F"="2.32"; "G"="2.42" },
    @{ "A"="23-08-2021"; "C"="-1.25"; "D"="0.68"; "E"="1.67"; "F"="2.35" },
    @{ "A"="24-08-2021"; "C"="-1.2"; "D"="0.67"; "E"="1.67"; "F"="2.33" },
    @{ "A"="25-08-2021"; "C"="-1.24"; "D"="0.63"; "E"="1.6" },
    @{ "A"="26-08-2021"; "C"="-1.28"; "D"="0.61"; "E"="1.6" },
    @{ "A"="27-08-2021"; "C"="-1.29"; "D"="0.57"; "E"="1.51"; "F"="2.28"; "G"="2.39" },
    @{ "A"="30-08-2021"; "C"="-1.26"; "D"="0.5"; "E"="1.45"; "F"="2.25" },
    @{ "A"="31-08-2021"; "C"="-1.32"; "D"="0.52"; "E"="1.46"; "F"="2.25"; "G"="2.36" },
    @{ "A"="01-09-2021"; "C"="-1.1"; "D"="0.71"; "E"="1.72"; "F"="2.4"; "G"="2.6" },
    @{ "A"="02-09-2021"; "C"="-0.83"; "D"="0.95"; "E"="1.83"; "F"="2.48" }
)

$startRow = 148
$endRow = $startRow + $newData.Count - 1
$cols = @("A","B","C","D","E","F","G")

# Column A holds dates formatted as plain text (e.g. "04-08-2021"), matching the
# existing rows above. Mark the range as Text first so Excel doesn't silently
# reinterpret values such as "04-08-2021" as a real date while we type them in.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowNum = $startRow + $i
    $entry = $newData[$i]
    foreach ($col in $cols) {
        if ($entry.ContainsKey($col)) {
            $value = $entry[$col]
            if ($col -eq "A") {
                $ws.Cells.Item($rowNum, 1).Value = [string]$value
            } else {
                $ws.Cells.Item($rowNum, [int][char]$col - [int][char]"A" + 1).Value = [double]$value
            }
        }
    }
}

# Put the cell style back to the workbook default so the new rows keep the same
# unstyled look as the rest of the data (only the header row uses a named style).
$ws.Range("A$startRow`:A$endRow").Style = "Normal"
